$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the team members in column B (rows 12-16)
$ws.Range("B12").Value = "Veselin"
$ws.Range("B13").Value = "Rawda"
$ws.Range("B14").Value = "Hannah"
$ws.Range("B15").Value = "Mirit"
$ws.Range("B16").Value = "Bogdana"

# Add new row 17 for a new team member
$ws.Range("B17").Value = "Martin"
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 1

# Explicitly set font color to black for the edited name cells (B12:B17)
$ws.Range("B12:B17").Font.Color = 0

# Update selection to match the diff (activeCell F13, sqref F13)
$ws.Range("F13").Select() | Out-Null
